$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data per GitHub Actions scrape run.
# D-column cells may look numeric (e.g. "1.000", "5.289"); force them to stay
# text so trailing zeros / grouping dots survive, then restore default style
# so no stray formatting is introduced.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.930.67"
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = "  -2.72%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.860.32"
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = "  -2.13%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("E5").Value = "  -2.05%  "

$ws.Range("E6").Value = "  +0.04%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5040"
$ws.Range("D7").Style = "Normal"

$ws.Range("E7").Value = "  -3.30%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3736"
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").Value = "  -1.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07131"
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").Value = "  -1.49%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8827"
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = "  -1.34%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.57"
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").Value = "  -2.44%  "

$ws.Range("E12").Value = "  -0.87%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.859.09"
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Value = "  -2.39%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.289"
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = "  -2.74%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.02"
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = "  -3.10%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.000"
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = "  +0.02%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008398"
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = "  -3.58%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.08"
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").Value = "  -2.59%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = "  +0.09%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.975.82"
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").Value = "  -2.68%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.027"
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = "  -1.94%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.101.19"
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").Value = "  -1.66%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.46"
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").Value = "  -3.31%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.448"
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").Value = "  -1.89%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.840"
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").Value = "  -1.35%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "147.01"
$ws.Range("D26").Style = "Normal"

$ws.Range("E26").Value = "  -3.96%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.92"
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").Value = "  -1.89%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.098"
$ws.Range("D28").Style = "Normal"

$ws.Range("E28").Value = "  -2.79%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "112.50"
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").Value = "  -1.73%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.661"
$ws.Range("D30").Style = "Normal"

$ws.Range("E30").Value = "  -3.49%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.680"
$ws.Range("D31").Style = "Normal"

$ws.Range("E31").Value = "  -2.68%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09033"
$ws.Range("D32").Style = "Normal"

$ws.Range("E32").Value = "  +0.54%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05127"
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").Value = "  -2.77%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.019"
$ws.Range("D34").Style = "Normal"

$ws.Range("E34").Value = "  -4.78%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.148"
$ws.Range("D35").Style = "Normal"

$ws.Range("E35").Value = "  -7.33%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7223"
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").Value = "  -6.59%  "

$ws.Range("E37").Value = "  -1.99%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.037"
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").Value = "  -0.47%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.451"
$ws.Range("D39").Style = "Normal"

$ws.Range("E39").Value = "  -6.06%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.075"
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").Value = "  -1.40%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5273"
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").Value = "  -3.95%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.508"
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").Value = "  -1.96%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "114.76"
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Value = "  +1.38%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.233"
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Value = "  -2.48%  "

$ws.Range("E45").Value = "  -2.61%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9995"
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").Value = "  +0.04%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4599"
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").Value = "  -3.80%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.938"
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").Value = "  -4.07%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.563"
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").Value = "  -3.05%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.47"
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").Value = "  -0.74%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.80"
$ws.Range("D51").Style = "Normal"

$ws.Range("E51").Value = "  -4.07%  "
